$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.257.14'
$ws.Range('E2').Value = '  +2.12%  '
$ws.Range('D3').Value = '2.364.45'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  +0.00%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '0.677'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +4.30%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '239.15'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +2.89%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '74.13'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +10.07%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +21.01%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.101'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +5.69%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '30.72'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +16.88%  '
$ws.Range('E12').Value = '  +2.37%  '
$ws.Range('D13').Value = '2.711.28'
$ws.Range('E13').Value = '  +0.34%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '16.91'
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '6.92'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +10.67%  '
$ws.Range('E16').Value = '  +8.10%  '
$ws.Range('D17').Value = '2.364.62'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').Value = '44.290.42'
$ws.Range('E18').Value = '  +2.35%  '
$ws.Range('E19').Value = '  +4.58%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '77.61'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +5.00%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.51'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +4.22%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '255.08'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +2.51%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '3.84'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -3.43%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E25').Value = '  +3.03%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '10.39'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +5.10%  '
$ws.Range('E27').Value = '  +3.86%  '
$ws.Range('E28').Value = '  +1.24%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '173.94'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.15%  '
$ws.Range('E30').Value = '  +3.56%  '
$ws.Range('E31').Value = '  +3.76%  '
$ws.Range('E32').Value = '  +5.43%  '
$ws.Range('E33').Value = '  +7.48%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '5.22'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +3.80%  '
$ws.Range('E35').Value = '  +2.98%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '3.91'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +8.06%  '
$ws.Range('E37').Value = '  -1.81%  '
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('E39').Value = '  +6.47%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '19.27'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +5.65%  '
$ws.Range('E41').Value = '  +0.17%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '8.86'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -1.16%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.26'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +3.17%  '
$ws.Range('E44').Value = '  +4.37%  '
$ws.Range('E45').Value = '  +1.03%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.187'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +12.75%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '99.41'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.89%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '4.47'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('E49').Value = '  +5.38%  '
$ws.Range('D50').Value = '1.447.07'
$ws.Range('D51').Value = '2.585.36'
$ws.Range('E51').Value = '  +0.41%  '
